# Generate Report for Handback
# The f79a32b9-0e45-419f-b9f3-90a59f4ef055.md file has now been handed back
# (it was previously shown as "Ready for handoff" with a stale-handback
# error). Update the status / handback datetime / error detail for that
# row across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K3").Value = "2016-08-17 12:45:54"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 13.67

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("K3").Value = "2016-08-17 12:46:06"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).ColumnWidth = 13.67
